# feat: add 2022-Q4 data
#
# 1) "总计" sheet: insert the 2022-Q4 summary as the new row 2, pushing the
#    existing 2022-Q2 summary row down to row 3.
# 2) Insert a brand new "2022-Q4" worksheet (positioned between "总计" and
#    "2022-Q2") holding the per-fund detail rows for the quarter.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------------
# 1) Update "总计": shift the old 2022-Q2 row to row 3, write 2022-Q4 on row 2
# ---------------------------------------------------------------------------

# Row 3 <- previous row 2 content (2022-Q2 summary)
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2022-Q2"
$wsTotal.Cells.Item(3, 3).Value = 5
$wsTotal.Cells.Item(3, 4).Value = 0.02

# Row 2 <- new 2022-Q4 summary
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 2
$wsTotal.Cells.Item(2, 4).Value = 0.01

# Carry the index-column style (bold/bordered, same as A2) onto the new A3
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right before "2022-Q2"
# ---------------------------------------------------------------------------

$newSheet = $wb.Worksheets.Add($wsQ2)
$newSheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
  $newSheet.Cells.Item(1, 2 + $i).Value = "'" + $headers[$i]
}

$fundRows = @(
  @(0, "009658", "汇丰晋信中小盘低波动策略股票A", "0.85", "92.42", "1.42", "0.0121", 10),
  @(1, "009775", "汇丰晋信中小盘低波动策略股票C", "0.04", "92.42", "1.42", "0.0006", 10)
)

for ($r = 0; $r -lt $fundRows.Length; $r++) {
  $rowData = $fundRows[$r]
  $excelRow = 2 + $r
  $newSheet.Cells.Item($excelRow, 1).Value = $rowData[0]
  $newSheet.Cells.Item($excelRow, 2).Value = "'" + $rowData[1]
  $newSheet.Cells.Item($excelRow, 3).Value = "'" + $rowData[2]
  $newSheet.Cells.Item($excelRow, 4).Value = "'" + $rowData[3]
  $newSheet.Cells.Item($excelRow, 5).Value = "'" + $rowData[4]
  $newSheet.Cells.Item($excelRow, 6).Value = "'" + $rowData[5]
  $newSheet.Cells.Item($excelRow, 7).Value = "'" + $rowData[6]
  $newSheet.Cells.Item($excelRow, 8).Value = $rowData[7]
}

# Header (B1:H1) gets the same bold/bordered style used by "总计"'s header row
$wsTotal.Range("B1:D1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Index column (A2:A3) gets the same style as "总计"!A2
$wsTotal.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# Remaining data cells (B2:H3) keep the plain/default style, like "总计"!C2
$wsTotal.Range("C2").Copy()
$newSheet.Range("B2:H3").PasteSpecial(-4122)

[void]$newSheet.Range("A1").Select()
